$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Superstore / Food, 2019-09-10): the rice cooker purchased here was later
# returned. Fold the refund into this row's cost instead of keeping a separate
# "Cookware" / "Returned the rice cooker" line further down, and update the
# comment to reflect that the item was returned.
$ws.Range("D7").Formula = "=-287.52+68.23"
$ws.Range("F7").Value = "Bough rice cooker, then returned later."

# Remove the now-redundant row that recorded the return of the rice cooker
# (2019-09-16, Superstore, Cookware, +68.23, "Returned the rice cooker").
# All rows below shift up by one.
$ws.Rows("28:28").Delete()

# Restore the selection reported in the saved file.
$ws.Range("F8").Select()
